# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the regenerated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value  = 1284
$wsExhibition.Range("F5").Value  = 1149
$wsExhibition.Range("F6").Value  = 14271
$wsExhibition.Range("F7").Value  = 16241
$wsExhibition.Range("F20").Value = 1242
$wsExhibition.Range("F23").Value = 29
$wsExhibition.Range("F24").Value = 6513
$wsExhibition.Range("F25").Value = 967
$wsExhibition.Range("F28").Value = 5
$wsExhibition.Range("F29").Value = 5695
$wsExhibition.Range("F30").Value = 97
$wsExhibition.Range("F32").Value = 172
$wsExhibition.Range("F33").Value = 4724

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 1284
$wsAll.Range("F5").Value  = 1149
$wsAll.Range("F6").Value  = 14271
$wsAll.Range("F7").Value  = 16241
$wsAll.Range("F20").Value = 1242
$wsAll.Range("F24").Value = 29
$wsAll.Range("F25").Value = 6513
$wsAll.Range("F26").Value = 967
$wsAll.Range("F29").Value = 5
$wsAll.Range("F31").Value = 5695
$wsAll.Range("F32").Value = 97
$wsAll.Range("F34").Value = 172
$wsAll.Range("F35").Value = 4724
